$d = $word.ActiveDocument

# Locate the paragraph that contains the "Gott begleitet ..." sentence.
# It currently sits between two empty paragraphs; the edit removes those
# two empty paragraph marks so the three paragraphs collapse into one.
$targetIndex = 0
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text -like "*Gott begleitet uns in unserer Trauer*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -gt 0) {
    $paragraphs = $d.Paragraphs

    $nextPara = $paragraphs.Item($targetIndex + 1)
    if ($nextPara.Range.Text.Trim().Length -eq 0) {
        $nextPara.Range.Delete()
    }

    $prevPara = $paragraphs.Item($targetIndex - 1)
    if ($prevPara.Range.Text.Trim().Length -eq 0) {
        $prevPara.Range.Delete()
    }
}
